$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1347.7142
$ws.Range("J17").Value = 1347.7142
$ws.Range("L17").Value = 4043.1426
$ws.Range("N17").Value = -4379.142599999999
$ws.Range("H43").Value = 13519.667
$ws.Range("I43").Value = 11000
$ws.Range("J43").Value = 13834.625
$ws.Range("K43").Value = 11000
$ws.Range("L43").Value = 13834.625
$ws.Range("M43").Value = -10931
$ws.Range("N43").Value = -13972.625
$ws.Range("H99").Value = 1641008.9
$ws.Range("J99").Value = 1289.6666
$ws.Range("L99").Value = 3868.9998
$ws.Range("N99").Value = -6864.9998
$ws.Range("H132").Value = 3034212.2
$ws.Range("I132").Value = 4025.2068
$ws.Range("K132").Value = 12075.6204
$ws.Range("M132").Value = -9545.6204
$ws.Range("H137").Value = 10332.042
$ws.Range("I137").Value = 12262.632
$ws.Range("K137").Value = 36787.896
$ws.Range("M137").Value = -34237.896
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 8553.983
$ws.Range("I32").Value = 8534.437
$ws.Range("J32").Value = 8733.166999999999
$ws.Range("K32").Value = 8534.437
$ws.Range("L32").Value = 8733.166999999999
$ws.Range("M32").Value = -8247.437
$ws.Range("N32").Value = -9307.166999999999
$ws.Range("H61").Value = 9643.875
$ws.Range("I61").Value = 10866.315
$ws.Range("K61").Value = 10866.315
$ws.Range("M61").Value = -10654.315
$ws.Range("H74").Value = 3881.075
$ws.Range("I74").Value = 4809.75
$ws.Range("K74").Value = 4809.75
$ws.Range("M74").Value = -3935.75
$ws.Range("H77").Value = 3881.075
$ws.Range("I77").Value = 4809.75
$ws.Range("K77").Value = 24048.75
$ws.Range("M77").Value = -19680.75
$ws.Range("H102").Value = 6624.6665
$ws.Range("I102").Value = 6725.1313
$ws.Range("K102").Value = 6725.1313
$ws.Range("M102").Value = -5103.1313
$ws.Range("H122").Value = 1310096.5
$ws.Range("I122").Value = 5511.5
$ws.Range("K122").Value = 16534.5
$ws.Range("M122").Value = -14084.5
$ws.Range("H136").Value = 9643.875
$ws.Range("I136").Value = 10866.315
$ws.Range("K136").Value = 32598.945
$ws.Range("M136").Value = -30048.945
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H130").Value = 73832.664
$ws.Range("J130").Value = 73832.664
$ws.Range("L130").Value = 73832.664
$ws.Range("N130").Value = -83872.664
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 8877.632
$ws.Range("I31").Value = 9229.75
$ws.Range("K31").Value = 9229.75
$ws.Range("M31").Value = -8934.75
$ws.Range("H34").Value = 8877.632
$ws.Range("I34").Value = 9229.75
$ws.Range("K34").Value = 9229.75
$ws.Range("M34").Value = -9027.75
$ws.Range("H99").Value = 360051.56
$ws.Range("J99").Value = 4800
$ws.Range("L99").Value = 4800
$ws.Range("N99").Value = -7796
$ws.Range("H105").Value = 402210
$ws.Range("I105").Value = 402210
$ws.Range("K105").Value = 402210
$ws.Range("M105").Value = -400463
$ws.Range("H122").Value = 5624.08
$ws.Range("I122").Value = 8175.5625
$ws.Range("K122").Value = 24526.6875
$ws.Range("M122").Value = -22076.6875
$ws.Range("H126").Value = 360051.56
$ws.Range("J126").Value = 4800
$ws.Range("L126").Value = 14400
$ws.Range("N126").Value = -19340
$ws.Range("H138").Value = 69666
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 69666
$ws.Range("K138").Value = 0
$ws.Range("L138").Value = 69666
$ws.Range("M138").Value = $null
$ws.Range("N138").Value = -79946
$ws.Range("H141").Value = 281922.7
$ws.Range("J141").Value = 306580.78
$ws.Range("L141").Value = 306580.78
$ws.Range("N141").Value = -316940.78
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 51369.1
$ws.Range("I11").Value = 31332.031
$ws.Range("K11").Value = 93996.09299999999
$ws.Range("M11").Value = -93856.09299999999
$ws.Range("H12").Value = 58.909092
$ws.Range("I12").Value = 67.666664
$ws.Range("J12").Value = 55.625
$ws.Range("K12").Value = 202.999992
$ws.Range("L12").Value = 166.875
$ws.Range("M12").Value = -29.99999199999999
$ws.Range("N12").Value = -512.875
$ws.Range("H56").Value = 6021.393
$ws.Range("I56").Value = 6021.393
$ws.Range("K56").Value = 6021.393
$ws.Range("M56").Value = -5491.393
$ws.Range("H68").Value = 10301.934
$ws.Range("J68").Value = 12423.25
$ws.Range("L68").Value = 37269.75
$ws.Range("N68").Value = -38891.75
$ws.Range("H71").Value = 10301.934
$ws.Range("J71").Value = 12423.25
$ws.Range("L71").Value = 111809.25
$ws.Range("N71").Value = -119921.25
$ws.Range("H105").Value = 9974.368
$ws.Range("I105").Value = 9026
$ws.Range("K105").Value = 27078
$ws.Range("M105").Value = -24457
$ws.Range("H107").Value = 1681.3334
$ws.Range("J107").Value = 1681.3334
$ws.Range("L107").Value = 5044.0002
$ws.Range("N107").Value = -8884.0002
$ws.Range("H140").Value = 3199.7334
$ws.Range("I140").Value = 3071.1428
$ws.Range("K140").Value = 9213.428400000001
$ws.Range("M140").Value = -4033.428400000001
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H52").Value = 24701.75
$ws.Range("I52").Value = 20676.666
$ws.Range("K52").Value = 20676.666
$ws.Range("M52").Value = -20417.666
$ws.Range("H53").Value = 21666.666
$ws.Range("J53").Value = 25000
$ws.Range("L53").Value = 25000
$ws.Range("N53").Value = -26262
$ws.Range("H80").Value = 9241.362999999999
$ws.Range("I80").Value = 13037
$ws.Range("K80").Value = 13037
$ws.Range("M80").Value = -12039
$ws.Range("H83").Value = 9241.362999999999
$ws.Range("I83").Value = 13037
$ws.Range("K83").Value = 65185
$ws.Range("M83").Value = -60193
$ws.Range("H97").Value = 8871.799999999999
$ws.Range("J97").Value = 4749
$ws.Range("L97").Value = 4749
$ws.Range("N97").Value = -5741
$ws.Range("H102").Value = 24268.5
$ws.Range("I102").Value = 42904
$ws.Range("J102").Value = 5633
$ws.Range("K102").Value = 42904
$ws.Range("L102").Value = 5633
$ws.Range("M102").Value = -41282
$ws.Range("N102").Value = -8877
$ws.Range("H122").Value = 6725.7646
$ws.Range("I122").Value = 4723.3794
$ws.Range("K122").Value = 14170.1382
$ws.Range("M122").Value = -11720.1382
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 4180.8423
$ws.Range("I61").Value = 1110.7693
$ws.Range("K61").Value = 1110.7693
$ws.Range("M61").Value = -908.7692999999999
$ws.Range("H113").Value = 4180.8423
$ws.Range("I113").Value = 1110.7693
$ws.Range("K113").Value = 1110.7693
$ws.Range("M113").Value = 1059.2307
$ws.Range("H122").Value = 6039
$ws.Range("I122").Value = 6943.364
$ws.Range("K122").Value = 20830.092
$ws.Range("M122").Value = -18380.092
$ws.Range("H130").Value = 0
$ws.Range("J130").Value = 0
$ws.Range("L130").Value = 0
$ws.Range("N130").Value = $null
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 1247.6285
$ws.Range("J113").Value = 3133.75
$ws.Range("L113").Value = 9401.25
$ws.Range("N113").Value = -13741.25
$ws.Range("H122").Value = 6661.737
$ws.Range("I122").Value = 5123.125
$ws.Range("J122").Value = 7780.727
$ws.Range("K122").Value = 15369.375
$ws.Range("L122").Value = 23342.181
$ws.Range("M122").Value = -12919.375
$ws.Range("N122").Value = -28242.181
$ws.Range("H136").Value = 7958.125
$ws.Range("I136").Value = 5725.423
$ws.Range("J136").Value = 17633.166
$ws.Range("K136").Value = 17176.269
$ws.Range("L136").Value = 52899.49800000001
$ws.Range("M136").Value = -14626.269
$ws.Range("N136").Value = -57999.49800000001
